$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing "*" from the row-label strings in column A (rows 2-6)
$ws.Range("A2").Value = "Aid perpetuates poverty as it makes people`nfeel less responsible for themselves"
$ws.Range("A3").Value = "Aid is not effective as most of it is diverted"
$ws.Range("A4").Value = "Aid is a pressure tactic for high-income countries that`nprevents low-income countries from developing freely"
$ws.Range("A5").Value = "[Country] is not responsible for`nwhat happens in other countries"
$ws.Range("A6").Value = "Charity begins at home: there is already a lot`nto do to support the [country] people in need"

# Overwrite the figures in column B (rows 2-6) with the corrected precision values
$ws.Range("B2").Value = 0.29024310987286
$ws.Range("B3").Value = 0.40171049482437
$ws.Range("B4").Value = 0.156571321939449
$ws.Range("B5").Value = 0.44955113604178
$ws.Range("B6").Value = 0.629727265520621
